# "Set Excel Table defaults"
#
# Renames the "n" option (shorthand for "number of records") to the
# more descriptive "number" wherever it appears as a SampleMetric
# option/default, and refreshes the sample StartDate/EndDate on the
# Setup sheet.

$wb = $excel.ActiveWorkbook

# --- Setup sheet: bump the example Start/End dates -------------------
$wsSetup = $wb.Worksheets.Item("Setup")
$wsSetup.Range("E3").Value = "2023-03-01"
$wsSetup.Range("F3").Value = "2023-06-01"

# --- ValidationReviewInformation sheet --------------------------------
$wsVRI = $wb.Worksheets.Item("ValidationReviewInformation")

# C1 holds the instructional note for the "SampleMetric" column; it
# calls out the two valid option values in bold ("n" / "proportion").
# Rename the "n" option to "number" while keeping the same rich-text
# emphasis (bold "Options:-", bold option names).
$c1 = $wsVRI.Range("C1")
$c1.Value = "[Review sample metric. ]`n`nOptions:`n- number (for number of records)`n- proportion (% of records)"
$c1.Characters(27,10).Font.Bold = $true
$c1.Characters(38,6).Font.Bold = $true
$c1.Characters(71,10).Font.Bold = $true

# C3 is the actual default/example SampleMetric value for reviewer #1.
$wsVRI.Range("C3").Value = "number"

# Reflect where the user last left their selection on this sheet, then
# restore focus back to the Setup tab (the workbook's active sheet).
$wsVRI.Range("D8").Select() | Out-Null
$wsSetup.Activate() | Out-Null
